$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1990291262135922
$ws.Range("C2").Value = 0.5388349514563107
$ws.Range("J2").Value = 0.009708737864077669
$ws.Range("P2").Value = 0.116504854368932
$ws.Range("S2").Value = 0.1359223300970874

# Row 3
$ws.Range("B3").Value = 0.01666666666666667
$ws.Range("C3").Value = 0.03333333333333333
$ws.Range("J3").Value = 0.04166666666666666
$ws.Range("P3").Value = 0.7166666666666667
$ws.Range("S3").Value = 0.1916666666666667

# Row 6
$ws.Range("B6").Value = 0.00909090909090909
$ws.Range("D6").Value = 0.01363636363636364
$ws.Range("F6").Value = 0.04545454545454546
$ws.Range("J6").Value = 0.2909090909090909
$ws.Range("O6").Value = 0.02272727272727273
$ws.Range("Q6").Value = 0.1318181818181818
$ws.Range("R6").Value = 0.07727272727272727
$ws.Range("S6").Value = 0.4090909090909091

# Row 7
$ws.Range("B7").Value = 0.0892018779342723
$ws.Range("D7").Value = 0.03755868544600939
$ws.Range("F7").Value = 0.05633802816901409
$ws.Range("J7").Value = 0.1643192488262911
$ws.Range("O7").Value = 0.01408450704225352
$ws.Range("Q7").Value = 0.1502347417840376
$ws.Range("R7").Value = 0.05633802816901409
$ws.Range("S7").Value = 0.431924882629108

# Row 8
$ws.Range("B8").Value = 0.06991525423728813
$ws.Range("D8").Value = 0.01906779661016949
$ws.Range("E8").Value = 0.00211864406779661
$ws.Range("F8").Value = 0.05932203389830509
$ws.Range("J8").Value = 0.06991525423728813
$ws.Range("O8").Value = 0.01694915254237288
$ws.Range("Q8").Value = 0.1949152542372881
$ws.Range("R8").Value = 0.1355932203389831
$ws.Range("S8").Value = 0.4322033898305085

# Row 9
$ws.Range("B9").Value = 0.0625
$ws.Range("D9").Value = 0.01442307692307692
$ws.Range("F9").Value = 0.0673076923076923
$ws.Range("J9").Value = 0.1201923076923077
$ws.Range("O9").Value = 0.01442307692307692
$ws.Range("Q9").Value = 0.1826923076923077
$ws.Range("R9").Value = 0.1009615384615385
$ws.Range("S9").Value = 0.4375

# Row 10
$ws.Range("B10").Value = 0.08503100088573959
$ws.Range("D10").Value = 0.02037201062887511
$ws.Range("F10").Value = 0.06908768821966342
$ws.Range("J10").Value = 0.08768821966341896
$ws.Range("O10").Value = 0.01771479185119575
$ws.Range("Q10").Value = 0.2276350752878654
$ws.Range("R10").Value = 0.09477413640389726
$ws.Range("S10").Value = 0.3976970770593445

# Row 11
$ws.Range("G11").Value = 0.155688622754491
$ws.Range("J11").Value = 0.09281437125748503
$ws.Range("K11").Value = 0.2095808383233533
$ws.Range("L11").Value = 0.5239520958083832
$ws.Range("S11").Value = 0.01796407185628742

# Row 12
$ws.Range("G12").Value = 0.7318435754189944
$ws.Range("J12").Value = 0.2067039106145251
$ws.Range("K12").Value = 0.0223463687150838
$ws.Range("L12").Value = 0.0223463687150838
$ws.Range("S12").Value = 0.01675977653631285

# Row 13
$ws.Range("G13").Value = 0.6470588235294118
$ws.Range("J13").Value = 0.3137254901960784
$ws.Range("S13").Value = 0.0392156862745098

# Row 14
$ws.Range("F14").Value = 0.2
$ws.Range("G14").Value = 0.6
$ws.Range("J14").Value = 0.2

# Row 15
$ws.Range("F15").Value = 0.03365384615384615
$ws.Range("H15").Value = 0.1730769230769231
$ws.Range("I15").Value = 0.08653846153846154
$ws.Range("J15").Value = 0.2932692307692308
$ws.Range("K15").Value = 0.08653846153846154
$ws.Range("N15").Value = 0.004807692307692308
$ws.Range("O15").Value = 0.04326923076923077
$ws.Range("S15").Value = 0.2788461538461539

# Row 16
$ws.Range("F16").Value = 0.05673758865248227
$ws.Range("H16").Value = 0.148936170212766
$ws.Range("I16").Value = 0.09219858156028368
$ws.Range("J16").Value = 0.3191489361702128
$ws.Range("K16").Value = 0.1134751773049645
$ws.Range("M16").Value = 0.02127659574468085
$ws.Range("O16").Value = 0.06382978723404255
$ws.Range("S16").Value = 0.1843971631205674

# Row 17
$ws.Range("F17").Value = 0.02947845804988662
$ws.Range("H17").Value = 0.2063492063492063
$ws.Range("I17").Value = 0.108843537414966
$ws.Range("J17").Value = 0.3922902494331066
$ws.Range("K17").Value = 0.09523809523809523
$ws.Range("M17").Value = 0.02040816326530612
$ws.Range("O17").Value = 0.04081632653061224
$ws.Range("S17").Value = 0.1065759637188209

# Row 18
$ws.Range("F18").Value = 0.01345291479820628
$ws.Range("H18").Value = 0.2331838565022422
$ws.Range("I18").Value = 0.07623318385650224
$ws.Range("J18").Value = 0.3497757847533632
$ws.Range("K18").Value = 0.07174887892376682
$ws.Range("M18").Value = 0.03587443946188341
$ws.Range("O18").Value = 0.09417040358744394
$ws.Range("S18").Value = 0.1255605381165919

# Row 19
$ws.Range("F19").Value = 0.01579778830963665
$ws.Range("H19").Value = 0.2187993680884676
$ws.Range("I19").Value = 0.09004739336492891
$ws.Range("J19").Value = 0.3420221169036335
$ws.Range("K19").Value = 0.1287519747235387
$ws.Range("M19").Value = 0.02606635071090047
$ws.Range("N19").Value = 0.00315955766192733
$ws.Range("O19").Value = 0.06872037914691943
$ws.Range("S19").Value = 0.1066350710900474
